$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.430.90"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "3.834.04"
$ws.Range("E3").Value = "  -2.48%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'511.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.07%  "

$ws.Range("D6").Value = "'138.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.36%  "

$ws.Range("E7").Value = "  -3.45%  "

$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("E9").Value = "  -5.11%  "

$ws.Range("E10").Value = "  -7.18%  "

$ws.Range("D11").Value = "'0.0000316"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.22%  "

$ws.Range("D12").Value = "'41.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.56%  "

$ws.Range("D13").Value = "'10.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D14").Value = "4.445.15"
$ws.Range("E14").Value = "  -2.42%  "

$ws.Range("D15").Value = "'21.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.06%  "

$ws.Range("D16").Value = "3.850.73"
$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("D17").Value = "'14.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "

$ws.Range("E18").Value = "  -2.14%  "

$ws.Range("E19").Value = "  +3.87%  "

$ws.Range("D20").Value = "68.411.18"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("D21").Value = "'414.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.57%  "

$ws.Range("D22").Value = "'3.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.90%  "

$ws.Range("D23").Value = "'13.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.57%  "

$ws.Range("D24").Value = "'85.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.20%  "

$ws.Range("D25").Value = "'3.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.63%  "

$ws.Range("D26").Value = "'11.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.31%  "

$ws.Range("D27").Value = "'10.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.71%  "

$ws.Range("D28").Value = "'35.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.49%  "

$ws.Range("D29").Value = "'672.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.59%  "

$ws.Range("D30").Value = "'13.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.67%  "

$ws.Range("E31").Value = "  -6.41%  "

$ws.Range("D32").Value = "'2.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("D33").Value = "'65.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.76%  "

$ws.Range("D34").Value = "'6.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.01%  "

$ws.Range("D35").Value = "'0.433"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.15%  "

$ws.Range("D36").Value = "'39.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.09%  "

$ws.Range("D37").Value = "0.0₃0819"
$ws.Range("E37").Value = "  -10.05%  "

$ws.Range("D38").Value = "'0.147"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").Value = "'3.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.48%  "

$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("E42").Value = "  -4.07%  "

$ws.Range("E43").Value = "  +2.75%  "

$ws.Range("D44").Value = "'2.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.35%  "

$ws.Range("E45").Value = "  -2.21%  "

$ws.Range("E46").Value = "  -3.90%  "

$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("D48").Value = "2.715.14"
$ws.Range("E48").Value = "  +11.74%  "

$ws.Range("D49").Value = "'143.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.72%  "

$ws.Range("E50").Value = "  -4.17%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'2.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.89%  "

Write-Host "Applied all changes"
